# TodoListe Projekt.xlsx - mark additional checklist items as done,
# remove a stray checkbox value, add a new (unchecked) checkbox entry,
# and move the current selection to reflect where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Validierungsregeln" block (rows 23-28): mark all checkboxes as done
$ws.Range("B23:B28").Value = $true

# "Repository-Layer" block: the heading row's checkbox (B64) should not
# carry a value at all (it belongs to the section title, not a task)
$ws.Range("B64").Clear()

# Mark the remaining repository tasks as done
$ws.Range("B66:B68").Value = $true

# New task row (Git & Versionskontrolle) gets an (unchecked) checkbox
$ws.Range("B70").Value = $false

# Reflect where the author ended up working in the sheet
$ws.Range("D68").Select()
